# Small implementation handling 403
# Adds the new observation row (A4/B4) to the forecast input table and
# widens the FORECAST.LINEAR ranges on row 25 to pick it up, then restores
# the author's last active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data point: date 2023-03-08 (serial 44993) with value 700.
# Copy A3's formatting (the date number format) onto A4 before writing the
# value so the new cell renders as a date just like the rest of column A.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").Value = 44993
$ws.Range("B4").Value = 700

# Extend the forecast formula so it regresses over rows 1-4 instead of 1-3.
$ws.Range("B25").Formula = "=_xlfn.FORECAST.LINEAR(A25,B1:B4,A1:A4)"

# Match the final selection recorded in the saved workbook.
$ws.Range("L17").Select()
